$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("D3").Value = 44508
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 667

# Row 4 updates
$ws.Range("D4").Value = 44525
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 8000
$ws.Range("P4").Value = 533
